# Update automatico via Actualizar 04-15-2021 12-34-33
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 15; $r++) {
    $ws.Range("D$r").Value = 44301.5236966881
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Range("D$r").Value = 44301.50230611111
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Range("D$r").Value = 44301.48091100695
}
